# PROS-6581 - CCRU - new KPI tables and POS 2019
#
# Populates the (previously empty) "Sheet5" worksheet with the new
# "PoS 2019 ..." / "Promo Visit" KPI lookup table (columns A, H and the
# CONCATENATE-based J helper formula), then restores "Initial INSERTS" as
# the active sheet/tab (it was "kpi_level_2" before the edit).

$wb = $excel.ActiveWorkbook

$kpiNames = @(
    "PoS 2019 - FT - CAP",
    "PoS 2019 - FT NS - CAP",
    "PoS 2019 - FT NS - REG",
    "PoS 2019 - FT - REG",
    "PoS 2019 - IC Canteen - EDU",
    "PoS 2019 - IC Canteen - OTH",
    "PoS 2019 - IC HoReCa BarTavernClub - CAP",
    "PoS 2019 - IC HoReCa BarTavernClub - REG",
    "PoS 2019 - IC HoReCa RestCafeTea - CAP",
    "PoS 2019 - IC HoReCa RestCafeTea - REG",
    "PoS 2019 - IC Petroleum - CAP",
    "PoS 2019 - IC Petroleum - REG",
    "PoS 2019 - IC QSR",
    "PoS 2019 - MT Conv Big - CAP",
    "PoS 2019 - MT Conv Big - REG",
    "PoS 2019 - MT Conv Small - CAP",
    "PoS 2019 - MT Conv Small - REG",
    "PoS 2019 - MT Hypermarket - CAP",
    "PoS 2019 - MT Hypermarket - REG",
    "PoS 2019 - MT Supermarket - CAP",
    "PoS 2019 - MT Supermarket - REG",
    "Promo Visit"
)

$sheet5 = $wb.Worksheets.Item("Sheet5")

for ($i = 0; $i -lt $kpiNames.Length; $i++) {
    $row = $i + 1
    $name = $kpiNames[$i]

    $sheet5.Cells.Item($row, 1).Value = $name

    if ($row -le 21) {
        $sheet5.Cells.Item($row, 8).Value = $name
        $sheet5.Cells.Item($row, 10).Formula = '=CONCATENATE("''",H' + $row + ',".xlsx'',")'
    }
}

$sheet5.Range("J1").Select()

# Restore "Initial INSERTS" as the active sheet/tab.
$wb.Worksheets.Item("Initial INSERTS").Activate()
